# Continued E-Ink display transition
# Rebuild the BOM table: drop E-Ink-transition-superseded parts (Relay Driver,
# old Arduino Nano -> MCU, old LCD/Pot parts, old Voltage Regulator / Relay /
# Inductor links), add the new parts, and compact the table up (old row 20
# "e-Ink Display" moves up to row 15; rows 17-20 become empty again).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Wipe everything first (hyperlinks are tracked independently of cell
# content/formatting in this model, so they must be removed explicitly) ---
$ws.Range("A1:C20").Hyperlinks.Delete()
$ws.Range("A1:C20").Clear()

# --- Header row ---
$ws.Range("A1").Value = "Description"
$ws.Range("B1").Value = "Qty"
$ws.Range("C1").Value = "Link"

# --- Row 2: Relay (new DS1E-M-DC3V part) ---
$ws.Range("A2").Value = "Relay"
$ws.Range("B2").Value = 3
$ws.Range("C2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("C2"), "https://www.digikey.com/en/products/detail/panasonic-electric-works/DS1E-M-DC3V/647259") | Out-Null

# --- Row 3: 2 Position Screw Terminal ---
$ws.Range("A3").Value = "2 Position Screw Terminal"
$ws.Range("B3").Value = 3
$ws.Range("C3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("C3"), "https://www.digikey.com/en/products/detail/phoenix-contact/5442206/4390264") | Out-Null

# --- Row 4: Voltage Regulator (new LM2596S-3-3 part; text only, no live link) ---
$ws.Range("A4").Value = "Voltage Regulator"
$ws.Range("B4").Value = 1
$ws.Range("C4").Style = "Hyperlink"
$ws.Range("C4").Value = "https://www.digikey.com/en/products/detail/texas-instruments/LM2596S-3-3/3701219"

# --- Row 5: Capicitor 220uF ---
$ws.Range("A5").Value = "Capicitor 220uF"
$ws.Range("B5").Value = 2
$ws.Range("C5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("C5"), "https://www.digikey.com/en/products/detail/nichicon/UVY1J221MPD1TD/4328548") | Out-Null

# --- Row 6: Schottky Diode ---
$ws.Range("A6").Value = "Schottky Diode"
$ws.Range("B6").Value = 1
$ws.Range("C6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("C6"), "https://www.digikey.com/en/products/detail/smc-diode-solutions/31DQ05TA/12142429") | Out-Null

# --- Row 7: Inductor 68 uH (new value; no link yet) ---
$ws.Range("A7").Value = "Inductor 68 uH"
$ws.Range("B7").Value = 1
$ws.Range("C7").Style = "Hyperlink"

# --- Row 8: 5mm LED ---
$ws.Range("A8").Value = "5mm LED"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "Sourced Locally"

# --- Row 9: 220R 1/4w THT Resistor ---
$ws.Range("A9").Value = "220R 1/4w THT Resistor"
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = "Sourced Locally"

# --- Row 10: 10K 1/4w THT Resistor ---
$ws.Range("A10").Value = "10K 1/4w THT Resistor"
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = "Sourced Locally"

# --- Row 11: Rectifier ---
$ws.Range("A11").Value = "Rectifier"
$ws.Range("B11").Value = 1
$ws.Range("C11").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("C11"), "https://www.digikey.com/en/products/detail/mdd/KBP307-L/14825026") | Out-Null

# --- Row 15: e-Ink Display (moved up from old row 20; plain text link) ---
$ws.Range("A15").Value = "e-Ink Display"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "https://www.adafruit.com/product/4777"

# --- Row 16: MCU (replaces "Arduino Nano") ---
$ws.Range("A16").Value = "MCU"
$ws.Range("B16").Value = 2

# --- Column C width ---
$ws.Range("C1").ColumnWidth = 116.1667

# --- Selection ---
$ws.Range("A17").Select()
